# Updated symbol list on Wed Dec 21 10:52:16 UTC 2022 with GitHub Actions
#
# Refreshes the crypto price/volume snapshot on the active sheet:
#  - updates Price (column D) for a batch of existing rows
#  - rotates rows 19-25 up one position (TigerCash/HotbitToken/BitKan/NitroEx/
#    LEO/BTSEToken/One each shift into the row above, with "One" landing at
#    row 19 and "BTSEToken" at row 25), refreshing Coin/Link/Price/Volume(1h)
#    for each moved row
#
# Column D holds prices as text (e.g. "249.32"), so values are written with a
# leading apostrophe to keep them as text instead of being re-interpreted as
# numbers by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price-only refreshes ---
$ws.Range("D2").Value  = "'249.33"
$ws.Range("D3").Value  = "'22.60"
$ws.Range("D4").Value  = "'5.403"
$ws.Range("D5").Value  = "'0.05683"
$ws.Range("D7").Value  = "'6.336"
$ws.Range("D8").Value  = "'0.8048"
$ws.Range("D9").Value  = "'0.9200"
$ws.Range("D10").Value = "'0.1401"
$ws.Range("D11").Value = "'0.07432"
$ws.Range("D12").Value = "'0.03136"
$ws.Range("D13").Value = "'0.03035"
$ws.Range("D14").Value = "'0.09381"
$ws.Range("D15").Value = "'3.910"
$ws.Range("D16").Value = "'0.001571"
$ws.Range("D17").Value = "'0.04810"

# --- Rows 19-25 rotate up one slot (new coin "One" enters at the top) ---
$ws.Range("B19").Value = "One"
$ws.Range("C19").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D19").Value = "'0.0005851"
$ws.Range("E19").Value = "18OneONE"

$ws.Range("B20").Value = "TigerCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D20").Value = "'0.006466"
$ws.Range("E20").Value = "19TigerCashTCH"

$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").Value = "'0.004993"
$ws.Range("E21").Value = "20HotbitTokenHTB"

$ws.Range("B22").Value = "BitKan"
$ws.Range("C22").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D22").Value = "'0.001006"
$ws.Range("E22").Value = "21BitKanKAN"

$ws.Range("B23").Value = "NitroEx"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D23").Value = "'0.0001500"
$ws.Range("E23").Value = "22NitroExNTX"

$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").Value = "'3.699"
$ws.Range("E24").Value = "23LEOLEO"

$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D25").Value = "'2.187"
$ws.Range("E25").Value = "24BTSETokenBTSE"

# --- Remaining price-only refreshes ---
$ws.Range("D26").Value = "'0.3253"
$ws.Range("D40").Value = "'0.04009"
$ws.Range("D41").Value = "'0.006920"
$ws.Range("D42").Value = "'0.1073"
$ws.Range("D43").Value = "'0.002710"
$ws.Range("D44").Value = "'0.007986"
$ws.Range("D45").Value = "'0.00005756"
$ws.Range("D47").Value = "'0.4991"
$ws.Range("D48").Value = "'0.2071"
